$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $r = $d.Content
    $ok = $r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output ("NOT FOUND: " + $find)
    }
    return $ok
}

# 1) Title: "Web Site" -> "Website"
Replace-Text "Web Site" "Website"

# 2) Generic "web site" -> "website" (4 occurrences across paragraphs)
Replace-Text "web site" "website"

# 3) "What" paragraph rewrite
Replace-Text "The USAFA Class of 1965 website is a site dedicated to all the members of the USAFA Class of 1965 who started with us as we entered the Cadet Wing in August 1961.  It aims to provide a central" "The USAFA Class of 1965 website is dedicated to the USAFA Class of 1965 members who entered the Cadet Wing in August 1961.  It provides a central"

# 4) Background paragraph rewrite (multiple targeted replacements, applied in document order)
Replace-Text "  The site was started in the summer of  2010 by Tom Browder and Bill Peavy" "  Tom Browder and Bill Peavy"

Replace-Text ", both of CS-24, as a way to contribute to  our 45" ", both of CS-24, started the site in the summer of 2010 as a way to contribute to our 45"

Replace-Text " class reunion.  Our immediate goal was to try to determine the fate of all our classmates who entered with us in CS-24 in the fall of 1961, and we were able to report “24" " class reunion.  Our immediate goal was to determine the fate of all our classmates who entered in CS-24 in August, 1961.  We were able to report “24"

Replace-Text "we have found out the fate of all but 156 classmates out of the approximately 756 men who entered with us that fall." "we have determined  the fate of all but 156 of our original 756 classmates.  (Note: Squadron representative are listed on the website.)"

# 5) Private area paragraph
Replace-Text ", for authenticated classmates, is available to show contact data that we have.  Authentication is via " ", for authenticated classmates, displays up-to-date contact data.  Authentication is via "

# 6) Pictures paragraph
Replace-Text "Pictures of us from the 1962 Polaris" "Pictures from the 1962 Polaris"

# 7) Apply distinct formatting (color + shading) to the newly-added "(Note: ...)" note in the Background paragraph,
#    matching the source's copy-pasted-from-web styling.
$noteFind = $d.Content
$noteOk = $noteFind.Find.Execute("Squadron representative are listed on the website.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output ("note found: " + $noteOk)
if ($noteOk) {
    $noteRange = $d.Range($noteFind.Start, $noteFind.End)
    Write-Output ("note text: " + $noteRange.Text)
    $noteRange.Font.Color = 2236962
}
